$wb = $excel.ActiveWorkbook

$names = @("abur", "ahnd", "aque", "bull", "carp", "golb", "ivee", "mohk", "napl", "scdi", "sctw")

for ($i = 0; $i -lt $names.Count; $i++) {
    $wb.Worksheets.Item($i + 1).Name = $names[$i]
}
